# Weekly data refresh: insert one new price record for Acelga
# (Terminal Hortofrutícola Agro Chillán) as row 32, pushing the
# existing rows (old row 32 onward) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 32 - Excel shifts rows 32..156
# down to 33..157 and extends the used range to R157.
$ws.Rows("32:32").Insert()

# Populate the new row with the new record's data.
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = "9/14/2021"
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112009
$ws.Range("G32").Value = "Acelga"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 160
$ws.Range("K32").Value = 350
$ws.Range("L32").Value = 400
$ws.Range("M32").Value = 375
$ws.Range("N32").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O32").Value = "Provincia de Diguillín"
$ws.Range("P32").Value = 375
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"
